$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 3054
$ws.Range("L3").Value = 3096
$ws.Range("L4").Value = 807
$ws.Range("L6").Value = 2758
$ws.Range("L7").Value = 9890

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L2").Value = 37
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 182
$ws.Range("L3").Value = 207
$ws.Range("L6").Value = 176
$ws.Range("L7").Value = 631

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L3").Value = 93
$ws.Range("L7").Value = 233

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 130
$ws.Range("L7").Value = 455

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L2").Value = 53
$ws.Range("L3").Value = 39
$ws.Range("L7").Value = 131

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 110
$ws.Range("L4").Value = 23
$ws.Range("L7").Value = 364

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 67
$ws.Range("L7").Value = 182

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L7").Value = 334
$ws.Range("L8").Value = 631
$ws.Range("L9").Value = 61
$ws.Range("L10").Value = 64
$ws.Range("L11").Value = 166
$ws.Range("L19").Value = 279
$ws.Range("L20").Value = 249
$ws.Range("L23").Value = 104
$ws.Range("L29").Value = 538
$ws.Range("L33").Value = 455
$ws.Range("L34").Value = 64
$ws.Range("L37").Value = 364
$ws.Range("L39").Value = 4
$ws.Range("L41").Value = 45
$ws.Range("L49").Value = 54
$ws.Range("L51").Value = 120
$ws.Range("L52").Value = 195
$ws.Range("L53").Value = 111
$ws.Range("L54").Value = 200
$ws.Range("L60").Value = 60
$ws.Range("L63").Value = 29
$ws.Range("L64").Value = 64
$ws.Range("L65").Value = 182
$ws.Range("L67").Value = 361
$ws.Range("L69").Value = 28
$ws.Range("L76").Value = 132
$ws.Range("L77").Value = 59
$ws.Range("L78").Value = 121
$ws.Range("L79").Value = 260
$ws.Range("L81").Value = 10
$ws.Range("L83").Value = 233
$ws.Range("L84").Value = 99
$ws.Range("L85").Value = 497
$ws.Range("L89").Value = 133
$ws.Range("L90").Value = 96
$ws.Range("L91").Value = 140
$ws.Range("L92").Value = 28
$ws.Range("L95").Value = 131
$ws.Range("L97").Value = 90
$ws.Range("L101").Value = 9890

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L3").Value = 133
$ws.Range("L7").Value = 361

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L2").Value = 36
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 99

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L3").Value = 44
$ws.Range("L6").Value = 98
$ws.Range("L7").Value = 200

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 164
$ws.Range("L3").Value = 202
$ws.Range("L7").Value = 538

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L3").Value = 83
$ws.Range("L7").Value = 279

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L2").Value = 25
$ws.Range("L3").Value = 25
$ws.Range("L6").Value = 62
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 93
$ws.Range("L6").Value = 91

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("L2").Value = 28
$ws.Range("L7").Value = 64

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L6").Value = 35
$ws.Range("L7").Value = 121

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L3").Value = 40
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L3").Value = 92
$ws.Range("L7").Value = 260

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 64

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 80
$ws.Range("L6").Value = 68
$ws.Range("L7").Value = 249

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 103
$ws.Range("L3").Value = 101
$ws.Range("L6").Value = 96
$ws.Range("L7").Value = 334

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("L2").Value = 17
$ws.Range("L3").Value = 18
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 64

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("L2").Value = 1
$ws.Range("L6").Value = 4

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 59
$ws.Range("L7").Value = 166

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L2").Value = 19
$ws.Range("L6").Value = 49
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("L2").Value = 11
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L3").Value = 44
$ws.Range("L4").Value = 6

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L2").Value = 42
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 133

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L2").Value = 35
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 96

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L3").Value = 36
$ws.Range("L7").Value = 120

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L3").Value = 23
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 145
$ws.Range("L3").Value = 204
$ws.Range("L6").Value = 100
$ws.Range("L7").Value = 497

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("L3").Value = 23
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L2").Value = 68
$ws.Range("L3").Value = 57
$ws.Range("L7").Value = 195

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 10
